$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each touched cell to Text format before assigning, then restore the
# default "Normal" style so number-like strings (e.g. "320.85") are kept as
# literal text instead of being auto-converted to numeric values.

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.711.51'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('E2').Style = 'Normal'

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.470.60'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('E3').Style = 'Normal'

# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('E4').Style = 'Normal'

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.45%  '
$ws.Range('E5').Style = 'Normal'

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '92.01'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.18%  '
$ws.Range('E6').Style = 'Normal'

# Row 7
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E7').Style = 'Normal'

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.508'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.40%  '
$ws.Range('E9').Style = 'Normal'

# Row 10
$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('B10').Style = 'Normal'
$ws.Range('C10').NumberFormat = '@'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('C10').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0856'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.78%  '
$ws.Range('E10').Style = 'Normal'

# Row 11
$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'Avalanche'
$ws.Range('B11').Style = 'Normal'
$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('C11').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '32.89'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.10%  '
$ws.Range('E11').Style = 'Normal'

# Row 12
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.85%  '
$ws.Range('E12').Style = 'Normal'

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.850.66'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('E13').Style = 'Normal'

# Row 14
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('E14').Style = 'Normal'

# Row 15
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.91%  '
$ws.Range('E15').Style = 'Normal'

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.465.79'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('E16').Style = 'Normal'

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.790'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.92%  '
$ws.Range('E17').Style = 'Normal'

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '41.636.73'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.18%  '
$ws.Range('E18').Style = 'Normal'

# Row 19
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.81%  '
$ws.Range('E19').Style = 'Normal'

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0940'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.51%  '
$ws.Range('E20').Style = 'Normal'

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.95'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.09%  '
$ws.Range('E21').Style = 'Normal'

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.21'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.81%  '
$ws.Range('E22').Style = 'Normal'

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '239.52'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E23').Style = 'Normal'

# Row 24
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.10%  '
$ws.Range('E24').Style = 'Normal'

# Row 25
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.70%  '
$ws.Range('E25').Style = 'Normal'

# Row 26
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('E26').Style = 'Normal'

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.80'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.42%  '
$ws.Range('E27').Style = 'Normal'

# Row 28
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.60%  '
$ws.Range('E28').Style = 'Normal'

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.71'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.06%  '
$ws.Range('E29').Style = 'Normal'

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.16'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.78%  '
$ws.Range('E30').Style = 'Normal'

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '155.60'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.32%  '
$ws.Range('E31').Style = 'Normal'

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.43'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.30%  '
$ws.Range('E32').Style = 'Normal'

# Row 33
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('E33').Style = 'Normal'

# Row 34
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'Hedera'
$ws.Range('B34').Style = 'Normal'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('C34').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0764'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.20%  '
$ws.Range('E34').Style = 'Normal'

# Row 35
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('B35').Style = 'Normal'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('C35').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.57'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.43%  '
$ws.Range('E35').Style = 'Normal'

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '16.99'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.95%  '
$ws.Range('E36').Style = 'Normal'

# Row 37
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.39%  '
$ws.Range('E37').Style = 'Normal'

# Row 39
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.76%  '
$ws.Range('E39').Style = 'Normal'

# Row 40
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('E40').Style = 'Normal'

# Row 41
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.14%  '
$ws.Range('E41').Style = 'Normal'

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.31'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -8.18%  '
$ws.Range('E42').Style = 'Normal'

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.000.55'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.55%  '
$ws.Range('E43').Style = 'Normal'

# Row 44
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.10%  '
$ws.Range('E44').Style = 'Normal'

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.59'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.34%  '
$ws.Range('E45').Style = 'Normal'

# Row 46
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.49%  '
$ws.Range('E46').Style = 'Normal'

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.51'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +5.11%  '
$ws.Range('E47').Style = 'Normal'

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.728.59'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.95%  '
$ws.Range('E48').Style = 'Normal'

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '75.77'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +4.18%  '
$ws.Range('E49').Style = 'Normal'

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '97.12'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.17%  '
$ws.Range('E50').Style = 'Normal'

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '66.91'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.23%  '
$ws.Range('E51').Style = 'Normal'
